$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the four new error code rows (412204-412207)
$ws.Range("A6").Value = 412204
$ws.Range("B6").Value = "The 'Set Subscription Attribute' VI only allows setting the 'DeliveryPolicy' and 'RawMessageDelivery' attributes."

$ws.Range("A7").Value = 412205
$ws.Range("B7").Value = "EndpointDisabled - The specified endpoint is disabled."

$ws.Range("A8").Value = 412206
$ws.Range("B8").Value = "ParameterValueInvalid - A request parameter does not comply with the associated constraints."

$ws.Range("A9").Value = 412207
$ws.Range("B9").Value = "PlatformApplicationDisabled - The specified platform application is disabled."

# Match formatting of the existing rows (style index 3: left-aligned)
$ws.Range("A6:A9").HorizontalAlignment = -4131

# Update the selected cell to match the diff
$ws.Range("B12").Select()
